$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Footer "last modified" date field: 7/16/2021 -> 7/21/2021
#    This literal text lives in the cached <a:t> of the datetimeFigureOut
#    field inside the "Date Placeholder" shape on the Slide Master and on
#    every Slide Layout. Find it by its current text and overwrite it.
# ---------------------------------------------------------------------------
function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "7/16/2021") {
                $sub = $tr.Characters(1, $tr.Length)
                $sub.Text = "7/21/2021"
            }
        }
    }
}

Update-DateShape $p.SlideMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateShape $layouts.Item($li).Shapes
}

# ---------------------------------------------------------------------------
# 2) Slide 7 "Traffic mirroring session" caption -> "Traffic mirror session"
#    (the other slides use the plural "Traffic mirroring sessions" caption
#    and are left untouched).
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
for ($i = 1; $i -le $s7.Shapes.Count; $i++) {
    $shp = $s7.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "Traffic mirroring session") {
            $idx = $tr.Text.IndexOf("mirroring ")
            $sub = $tr.Characters($idx + 1, 10)
            $sub.Text = "mirror "
        }
    }
}
